$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.538.67"

$ws.Range("D3").Value = "2.109.02"
$ws.Range("E3").Value = "  +4.93%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5266"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4366"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08886"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.49%  "

$ws.Range("E11").Value = "  +2.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.58"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "2.107.85"
$ws.Range("E13").Value = "  +4.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.732"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.769"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001129"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06639"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.322"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("D23").Value = "30.588.59"
$ws.Range("E23").Value = "  +0.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.51%  "

$ws.Range("D26").Value = "2.353.79"
$ws.Range("E26").Value = "  +4.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.598"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.204"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.84%  "

$ws.Range("E32").Value = "  +2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.677"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +22.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.210"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.926"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "

$ws.Range("E36").Value = "  +9.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02582"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.497"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.10%  "

$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2279"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6809"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.261"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6369"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.210"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.626"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.197"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.04%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "
